$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")
$ws.Rows.Item(4).Delete() | Out-Null
$ws.Range("B13").Select() | Out-Null
